$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2 (source row 2)
$ws.Range("C2").Value = 46073

# row 3 (source row 3)
$ws.Range("C3").Value = 46073

# row 4 (source row 4)
$ws.Range("C4").Value = 46073

# row 5 (source row 5)
$ws.Range("C5").Value = 46073

# row 6 (source row 6)
$ws.Range("C6").Value = 46073

# row 7 (source row 7)
$ws.Range("C7").Value = 46073

# row 8 (source row 8)
$ws.Range("C8").Value = 46073

# row 9 (source row 9)
$ws.Range("C9").Value = 46073

# row 10 (source row 10)
$ws.Range("C10").Value = 46073

# row 11 (source row 11)
$ws.Range("C11").Value = 46073

# row 12 (source row 12)
$ws.Range("C12").Value = 46073

# row 13 (source row 13)
$ws.Range("C13").Value = 46073

# row 14 (source row 42)
$ws.Range("A14").Value = "A 53343-2024"
$ws.Range("B14").Value = 45614.43885416666
$ws.Range("C14").Value = 46073
$ws.Range("F14").ClearContents()
$ws.Range("G14").Value = 0.9

# row 15 (source row 16)
$ws.Range("A15").Value = "A 10263-2024"
$ws.Range("B15").Value = 45365.43090277778
$ws.Range("C15").Value = 46073
$ws.Range("F15").Value = "Kyrkan"
$ws.Range("G15").Value = 1.4

# row 16 (source row 35)
$ws.Range("A16").Value = "A 30174-2021"
$ws.Range("B16").Value = 44363
$ws.Range("C16").Value = 46073
$ws.Range("F16").ClearContents()
$ws.Range("G16").Value = 1.8

# row 17 (source row 17)
$ws.Range("C17").Value = 46073

# row 18 (source row 18)
$ws.Range("C18").Value = 46073

# row 19 (source row 47)
$ws.Range("A19").Value = "A 34984-2024"
$ws.Range("B19").Value = 45527
$ws.Range("C19").Value = 46073
$ws.Range("G19").Value = 4.1

# row 20 (source row 22)
$ws.Range("A20").Value = "A 42994-2025"
$ws.Range("B20").Value = 45909.45351851852
$ws.Range("C20").Value = 46073
$ws.Range("G20").Value = 7.9

# row 21 (source row 21)
$ws.Range("C21").Value = 46073

# row 22 (source row 24)
$ws.Range("A22").Value = "A 43448-2025"
$ws.Range("B22").Value = 45911.45209490741
$ws.Range("C22").Value = 46073
$ws.Range("G22").Value = 1.2

# row 23 (source row 27)
$ws.Range("A23").Value = "A 44192-2025"
$ws.Range("B23").Value = 45915.61556712963
$ws.Range("C23").Value = 46073
$ws.Range("G23").Value = 0.8

# row 24 (source row 49)
$ws.Range("A24").Value = "A 22072-2023"
$ws.Range("B24").Value = 45069
$ws.Range("C24").Value = 46073
$ws.Range("G24").Value = 3.5

# row 25 (source row 23)
$ws.Range("A25").Value = "A 27365-2025"
$ws.Range("B25").Value = 45812.64355324074
$ws.Range("C25").Value = 46073
$ws.Range("G25").Value = 11.9

# row 26 (source row 20)
$ws.Range("A26").Value = "A 55562-2022"
$ws.Range("B26").Value = 44888
$ws.Range("C26").Value = 46073
$ws.Range("G26").Value = 0.8

# row 27 (source row 28)
$ws.Range("A27").Value = "A 58125-2025"
$ws.Range("B27").Value = 45982.61506944444
$ws.Range("C27").Value = 46073

# row 28 (source row 25)
$ws.Range("A28").Value = "A 58109-2025"
$ws.Range("B28").Value = 45982.59765046297
$ws.Range("C28").Value = 46073
$ws.Range("G28").Value = 2.5

# row 29 (source row 26)
$ws.Range("A29").Value = "A 58111-2025"
$ws.Range("B29").Value = 45982.59920138889
$ws.Range("C29").Value = 46073
$ws.Range("G29").Value = 0.6

# row 30 (source row 30)
$ws.Range("C30").Value = 46073

# row 31 (source row 32)
$ws.Range("A31").Value = "A 54203-2025"
$ws.Range("B31").Value = 45964
$ws.Range("C31").Value = 46073
$ws.Range("G31").Value = 10.3

# row 32 (source row 34)
$ws.Range("A32").Value = "A 33799-2025"
$ws.Range("B32").Value = 45842.44056712963
$ws.Range("C32").Value = 46073
$ws.Range("F32").Value = "Allmännings- och besparingsskogar"
$ws.Range("G32").Value = 2.5

# row 33 (source row 48)
$ws.Range("A33").Value = "A 16762-2022"
$ws.Range("B33").Value = 44673.47876157407
$ws.Range("C33").Value = 46073
$ws.Range("G33").Value = 4.2

# row 34 (source row 39)
$ws.Range("A34").Value = "A 1621-2026"
$ws.Range("B34").Value = 46034.47645833333
$ws.Range("C34").Value = 46073
$ws.Range("F34").ClearContents()
$ws.Range("G34").Value = 1.4

# row 35 (source row 40)
$ws.Range("A35").Value = "A 1622-2026"
$ws.Range("B35").Value = 46034.47929398148
$ws.Range("C35").Value = 46073
$ws.Range("G35").Value = 1.4

# row 36 (source row 36)
$ws.Range("C36").Value = 46073

# row 37 (source row 37)
$ws.Range("C37").Value = 46073

# row 38 (source row 29)
$ws.Range("A38").Value = "A 50230-2024"
$ws.Range("B38").Value = 45600
$ws.Range("C38").Value = 46073
$ws.Range("G38").Value = 1.7

# row 39 (source row 38)
$ws.Range("A39").Value = "A 13510-2025"
$ws.Range("B39").Value = 45736.47103009259
$ws.Range("C39").Value = 46073
$ws.Range("G39").Value = 2.5

# row 40 (source row 43)
$ws.Range("A40").Value = "A 62433-2025"
$ws.Range("B40").Value = 46007
$ws.Range("C40").Value = 46073
$ws.Range("G40").Value = 2

# row 41 (source row 15)
$ws.Range("A41").Value = "A 21972-2023"
$ws.Range("B41").Value = 45068.66232638889
$ws.Range("C41").Value = 46073
$ws.Range("F41").ClearContents()
$ws.Range("G41").Value = 1.5

# row 42 (source row 51)
$ws.Range("A42").Value = "A 35036-2024"
$ws.Range("B42").Value = 45527
$ws.Range("C42").Value = 46073
$ws.Range("G42").Value = 1.7

# row 43 (source row 52)
$ws.Range("A43").Value = "A 62831-2023"
$ws.Range("B43").Value = 45270
$ws.Range("C43").Value = 46073
$ws.Range("G43").Value = 1.1

# row 44 (source row 44)
$ws.Range("C44").Value = 46073

# row 45 (source row 50)
$ws.Range("A45").Value = "A 30766-2022"
$ws.Range("B45").Value = 44764
$ws.Range("C45").Value = 46073
$ws.Range("F45").ClearContents()
$ws.Range("G45").Value = 0.6

# row 46 (source row 46)
$ws.Range("C46").Value = 46073

# row 47 (source row 19)
$ws.Range("A47").Value = "A 30743-2021"
$ws.Range("B47").Value = 44365
$ws.Range("C47").Value = 46073
$ws.Range("G47").Value = 3

# row 48 (source row 45)
$ws.Range("A48").Value = "A 8436-2023"
$ws.Range("B48").Value = 44977
$ws.Range("C48").Value = 46073
$ws.Range("F48").Value = "Kyrkan"
$ws.Range("G48").Value = 4

# row 49 (source row 14)
$ws.Range("A49").Value = "A 21264-2022"
$ws.Range("B49").Value = 44705
$ws.Range("C49").Value = 46073
$ws.Range("F49").Value = "Allmännings- och besparingsskogar"
$ws.Range("G49").Value = 2.4

# row 50 (source row 31)
$ws.Range("A50").Value = "A 12077-2022"
$ws.Range("B50").Value = 44636.47484953704
$ws.Range("C50").Value = 46073
$ws.Range("G50").Value = 2.1

# row 51 (source row 33)
$ws.Range("A51").Value = "A 37407-2023"
$ws.Range("B51").Value = 45156.60152777778
$ws.Range("C51").Value = 46073
$ws.Range("G51").Value = 3.3

# row 52 (source row 41)
$ws.Range("A52").Value = "A 14149-2022"
$ws.Range("B52").Value = 44651
$ws.Range("C52").Value = 46073
$ws.Range("F52").Value = "Allmännings- och besparingsskogar"
$ws.Range("G52").Value = 3.8

